$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 'score' (column J) values for existing rows, per latest prediction run.
$scoreUpdates = @(
    @{ Row = 3; Value = 0 },
    @{ Row = 13; Value = 0 },
    @{ Row = 17; Value = 0.05333333333333334 },
    @{ Row = 19; Value = 0.006666666666666667 },
    @{ Row = 21; Value = 0.006666666666666667 },
    @{ Row = 29; Value = 0.02666666666666667 },
    @{ Row = 43; Value = 0.006666666666666667 },
    @{ Row = 62; Value = 0.006666666666666667 },
    @{ Row = 71; Value = 0 },
    @{ Row = 74; Value = 0.7066666666666667 },
    @{ Row = 81; Value = 0.04 },
    @{ Row = 84; Value = 0.07333333333333333 },
    @{ Row = 87; Value = 0 },
    @{ Row = 88; Value = 0.006666666666666667 },
    @{ Row = 96; Value = 0.8133333333333334 },
    @{ Row = 97; Value = 0 },
    @{ Row = 98; Value = 0.1933333333333333 },
    @{ Row = 102; Value = 0.01333333333333333 },
    @{ Row = 125; Value = 0.1733333333333333 },
    @{ Row = 126; Value = 0.006666666666666667 },
    @{ Row = 128; Value = 0.006666666666666667 },
    @{ Row = 132; Value = 0 },
    @{ Row = 138; Value = 0.006666666666666667 },
    @{ Row = 144; Value = 0.07333333333333333 },
    @{ Row = 145; Value = 0 },
    @{ Row = 151; Value = 0.006666666666666667 },
    @{ Row = 156; Value = 0.006666666666666667 },
    @{ Row = 173; Value = 0.76 },
    @{ Row = 175; Value = 0 },
    @{ Row = 191; Value = 0.05333333333333334 },
    @{ Row = 204; Value = 0.01333333333333333 },
    @{ Row = 207; Value = 0.006666666666666667 },
    @{ Row = 214; Value = 0 },
    @{ Row = 215; Value = 0.006666666666666667 },
    @{ Row = 216; Value = 0.006666666666666667 },
    @{ Row = 217; Value = 0 },
    @{ Row = 219; Value = 0 },
    @{ Row = 221; Value = 0.006666666666666667 },
    @{ Row = 232; Value = 0.02 },
    @{ Row = 235; Value = 0.06 },
    @{ Row = 237; Value = 0.04 },
    @{ Row = 238; Value = 0.03333333333333333 },
    @{ Row = 242; Value = 0.006666666666666667 },
    @{ Row = 247; Value = 0.5733333333333334 },
    @{ Row = 253; Value = 0 },
    @{ Row = 262; Value = 0.74 },
    @{ Row = 267; Value = 0.09333333333333334 },
    @{ Row = 271; Value = 0.03333333333333333 },
    @{ Row = 283; Value = 0.02 },
    @{ Row = 288; Value = 0 }
)

foreach ($u in $scoreUpdates) {
    $ws.Cells.Item($u.Row, 10).Value = $u.Value
}

# Add the new cryptocurrency row: Moo Deng
$newRow = 289
$ws.Cells.Item($newRow, 1).Value = "Moo Deng"
$ws.Cells.Item($newRow, 2).Value = "moodeng"
$ws.Cells.Item($newRow, 3).Value = "https://coin-images.coingecko.com/coins/images/50264/large/MOODENG.jpg?1726726975"
$ws.Cells.Item($newRow, 4).Value = "/charts/moodeng_chart.png"
$ws.Cells.Item($newRow, 5).Value = 188745454
$ws.Cells.Item($newRow, 6).Value = 0.189982
$ws.Cells.Item($newRow, 7).Value = 12.6837776520067
$ws.Cells.Item($newRow, 8).Value = 40.4236651718855
$ws.Cells.Item($newRow, 9).Value = -12.48789007620941
$ws.Cells.Item($newRow, 10).Value = 0.08
$ws.Cells.Item($newRow, 11).Value = $false
$ws.Cells.Item($newRow, 12).Value = ""
